# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 6 (Feria Lagunitas de Puerto Montt,
# Granada, "Primera"), shifting the existing rows 6-39 down to 7-40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 6 - this shifts rows 6..39 down to 7..40
# and keeps/propagates the date-column number format (style index 2) used by
# the rest of column D.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly record.
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C6").Value = "Los Lagos"
$ws.Range("D6").Value = 44659
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100104
$ws.Range("H6").Value = "Frutos de pepita"
$ws.Range("I6").Value = 100104001
$ws.Range("J6").Value = "Granada"
$ws.Range("K6").Value = "Wonderfull"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 14000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 14500
$ws.Range("Q6").Value = "$/caja 14 kilos empedrada"
$ws.Range("R6").Value = "Provincia de Limarí"
$ws.Range("S6").Value = 1036
$ws.Range("T6").Value = 14
